# Apply the updated cryptocurrency price/volume snapshot values.
# Numeric-looking text (e.g. "1.004", "0.00001091") is written with a leading
# apostrophe so Excel stores it verbatim as text (matching the original
# t="inlineStr" cells) instead of silently coercing it to a Double and losing
# formatting such as trailing zeros or the dotted "27.150.75" style values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.150.75'
$ws.Range("E2").Value = '  -1.21%  '

$ws.Range("D3").Value = '1.783.79'
$ws.Range("E3").Value = '  -1.93%  '

$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '''337.26'
$ws.Range("E5").Value = '  -1.88%  '

$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").Value = '''0.3881'
$ws.Range("E7").Value = '  +1.11%  '

$ws.Range("D8").Value = '''0.3436'
$ws.Range("E8").Value = '  -2.52%  '

$ws.Range("D9").Value = '''47.93'
$ws.Range("E9").Value = '  -2.09%  '

$ws.Range("D10").Value = '''1.188'
$ws.Range("E10").Value = '  -3.93%  '

$ws.Range("D11").Value = '''0.07446'
$ws.Range("E11").Value = '  -4.81%  '

$ws.Range("D12").Value = '''1.000'
$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("D13").Value = '''21.67'
$ws.Range("E13").Value = '  -2.96%  '

$ws.Range("D14").Value = '''6.430'
$ws.Range("E14").Value = '  -2.70%  '

$ws.Range("D15").Value = '1.783.15'
$ws.Range("E15").Value = '  -1.78%  '

$ws.Range("D16").Value = '''7.114'
$ws.Range("E16").Value = '  -1.67%  '

$ws.Range("D17").Value = '''0.00001091'
$ws.Range("E17").Value = '  -2.81%  '

$ws.Range("D18").Value = '''0.06646'
$ws.Range("E18").Value = '  -1.23%  '

$ws.Range("D19").Value = '''83.30'
$ws.Range("E19").Value = '  -3.75%  '

$ws.Range("D20").Value = '''1.000'
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("D21").Value = '''17.54'
$ws.Range("E21").Value = '  -0.76%  '

$ws.Range("D22").Value = '''6.511'
$ws.Range("E22").Value = '  -1.07%  '

$ws.Range("D23").Value = '27.163.84'

$ws.Range("E24").Value = '  -6.50%  '

$ws.Range("E25").Value = '  -3.77%  '

$ws.Range("D26").Value = '''21.16'
$ws.Range("E26").Value = '  -5.64%  '

$ws.Range("D27").Value = '''2.495'
$ws.Range("E27").Value = '  -7.17%  '

$ws.Range("D28").Value = '''1.444'
$ws.Range("E28").Value = '  -1.69%  '

$ws.Range("D29").Value = '''156.55'
$ws.Range("E29").Value = '  +1.64%  '

$ws.Range("D30").Value = '1.984.05'
$ws.Range("E30").Value = '  -1.76%  '

$ws.Range("D31").Value = '''134.02'
$ws.Range("E31").Value = '  -1.90%  '

$ws.Range("D32").Value = '''3.975'
$ws.Range("E32").Value = '  -2.15%  '

$ws.Range("D33").Value = '''5.979'
$ws.Range("E33").Value = '  -6.07%  '

$ws.Range("D34").Value = '''0.08697'
$ws.Range("E34").Value = '  -1.33%  '

$ws.Range("D35").Value = '''12.94'
$ws.Range("E35").Value = '  -7.03%  '

$ws.Range("D36").Value = '''1.620'
$ws.Range("E36").Value = '  -4.09%  '

$ws.Range("D37").Value = '''5.398'
$ws.Range("E37").Value = '  -4.52%  '

$ws.Range("D38").Value = '''0.6814'
$ws.Range("E38").Value = '  -2.84%  '

$ws.Range("D39").Value = '''0.06338'
$ws.Range("E39").Value = '  -2.55%  '

$ws.Range("D40").Value = '''0.02343'
$ws.Range("E40").Value = '  -2.63%  '

$ws.Range("E41").Value = '  -3.37%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''1.239'
$ws.Range("E42").Value = '  -4.26%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''8.434'
$ws.Range("E43").Value = '  -6.27%  '

$ws.Range("D44").Value = '''14.29'
$ws.Range("E44").Value = '  -2.94%  '

$ws.Range("D45").Value = '''0.9999'
$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("D46").Value = '''0.6398'
$ws.Range("E46").Value = '  -3.06%  '

$ws.Range("D47").Value = '''3.855'
$ws.Range("E47").Value = '  -2.64%  '

$ws.Range("D48").Value = '''2.148'
$ws.Range("E48").Value = '  -1.76%  '

$ws.Range("D49").Value = '''131.71'
$ws.Range("E49").Value = '  -0.87%  '

$ws.Range("D50").Value = '''0.07116'
$ws.Range("E50").Value = '  -2.89%  '

$ws.Range("D51").Value = '''79.23'
$ws.Range("E51").Value = '  -1.76%  '
